$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2252747252747253
$ws.Range("C2").Value = 0.4835164835164835
$ws.Range("J2").Value = 0.02197802197802198
$ws.Range("P2").Value = 0.1868131868131868
$ws.Range("S2").Value = 0.08241758241758242
# Row 3
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.04395604395604396
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.2032967032967033
# Row 4
$ws.Range("P4").Value = 0.6896551724137931
$ws.Range("S4").Value = 0.3103448275862069
# Row 6
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.005128205128205128
$ws.Range("F6").Value = 0.01538461538461539
$ws.Range("J6").Value = 0.3282051282051282
$ws.Range("O6").Value = 0.02051282051282051
$ws.Range("Q6").Value = 0.1692307692307692
$ws.Range("R6").Value = 0.04615384615384616
$ws.Range("S6").Value = 0.3487179487179487
# Row 7
$ws.Range("B7").Value = 0.1358695652173913
$ws.Range("D7").Value = 0.0108695652173913
$ws.Range("F7").Value = 0.03260869565217391
$ws.Range("J7").Value = 0.1467391304347826
$ws.Range("O7").Value = 0.01630434782608696
$ws.Range("Q7").Value = 0.1793478260869565
$ws.Range("R7").Value = 0.08152173913043478
$ws.Range("S7").Value = 0.3967391304347826
# Row 8
$ws.Range("B8").Value = 0.1100917431192661
$ws.Range("D8").Value = 0.02752293577981652
$ws.Range("F8").Value = 0.03669724770642202
$ws.Range("J8").Value = 0.1536697247706422
$ws.Range("O8").Value = 0.02752293577981652
$ws.Range("Q8").Value = 0.1674311926605505
$ws.Range("R8").Value = 0.0871559633027523
$ws.Range("S8").Value = 0.3899082568807339
# Row 9
$ws.Range("B9").Value = 0.1445086705202312
$ws.Range("D9").Value = 0.0115606936416185
$ws.Range("F9").Value = 0.06358381502890173
$ws.Range("J9").Value = 0.1502890173410405
$ws.Range("O9").Value = 0.04624277456647399
$ws.Range("Q9").Value = 0.1676300578034682
$ws.Range("R9").Value = 0.06936416184971098
$ws.Range("S9").Value = 0.3468208092485549
# Row 10
$ws.Range("B10").Value = 0.1198808637379002
$ws.Range("D10").Value = 0.03127326880119136
$ws.Range("E10").Value = 0.0007446016381236039
$ws.Range("F10").Value = 0.06031273268801191
$ws.Range("J10").Value = 0.1280714817572599
$ws.Range("O10").Value = 0.02382725241995532
$ws.Range("Q10").Value = 0.2196574832464631
$ws.Range("R10").Value = 0.08041697691734921
$ws.Range("S10").Value = 0.3358153387937454
# Row 11
$ws.Range("G11").Value = 0.130030959752322
$ws.Range("J11").Value = 0.1269349845201238
$ws.Range("K11").Value = 0.2136222910216718
$ws.Range("L11").Value = 0.4984520123839009
$ws.Range("S11").Value = 0.03095975232198142
# Row 12
$ws.Range("G12").Value = 0.7380952380952381
$ws.Range("J12").Value = 0.2023809523809524
$ws.Range("K12").Value = 0.005952380952380952
$ws.Range("L12").Value = 0.01785714285714286
$ws.Range("S12").Value = 0.03571428571428571
# Row 13
$ws.Range("G13").Value = 0.5384615384615384
$ws.Range("J13").Value = 0.3846153846153846
$ws.Range("S13").Value = 0.07692307692307693
# Row 15
$ws.Range("F15").Value = 0.03636363636363636
$ws.Range("H15").Value = 0.1318181818181818
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.3
$ws.Range("K15").Value = 0.05909090909090909
$ws.Range("M15").Value = 0.00909090909090909
$ws.Range("O15").Value = 0.05
$ws.Range("S15").Value = 0.3636363636363636
# Row 16
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1361702127659574
$ws.Range("I16").Value = 0.08936170212765958
$ws.Range("J16").Value = 0.4297872340425532
$ws.Range("K16").Value = 0.1063829787234043
$ws.Range("N16").Value = 0.008510638297872341
$ws.Range("O16").Value = 0.06382978723404255
$ws.Range("S16").Value = 0.1446808510638298
# Row 17
$ws.Range("F17").Value = 0.02637362637362637
$ws.Range("H17").Value = 0.1516483516483516
$ws.Range("I17").Value = 0.07692307692307693
$ws.Range("J17").Value = 0.4483516483516484
$ws.Range("K17").Value = 0.0989010989010989
$ws.Range("M17").Value = 0.01538461538461539
$ws.Range("O17").Value = 0.06593406593406594
$ws.Range("S17").Value = 0.1164835164835165
# Row 18
$ws.Range("F18").Value = 0.005555555555555556
$ws.Range("H18").Value = 0.1388888888888889
$ws.Range("I18").Value = 0.06111111111111111
$ws.Range("J18").Value = 0.3888888888888889
$ws.Range("K18").Value = 0.1333333333333333
$ws.Range("M18").Value = 0.02222222222222222
$ws.Range("O18").Value = 0.08333333333333333
$ws.Range("S18").Value = 0.1666666666666667
# Row 19
$ws.Range("F19").Value = 0.01795918367346939
$ws.Range("H19").Value = 0.2342857142857143
$ws.Range("I19").Value = 0.07836734693877551
$ws.Range("J19").Value = 0.3722448979591837
$ws.Range("K19").Value = 0.1126530612244898
$ws.Range("M19").Value = 0.02204081632653061
$ws.Range("O19").Value = 0.0546938775510204
$ws.Range("S19").Value = 0.1077551020408163
